$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Set new cell data for rows 2-8 (existing rows 2-4 updated, new rows 5-8 added)
# Code (col A) and Phone (col D) are stored as text in the source data, so we
# force text entry using a leading apostrophe, then clear the resulting
# "quote prefix" style so the cell keeps style index 0 like the rest.

$ws.Range("A2").Value = "'1"
$ws.Range("B2").Value = "sanjeev praj"
$ws.Range("C2").Value = "prajapatisanjiv8@gmail.com"
$ws.Range("D2").Value = "'8799879871"
$ws.Range("E2").Value = 100090

$ws.Range("A3").Value = "'2"
$ws.Range("B3").Value = "balasaheb"
$ws.Range("C3").Value = "balasaheb.more@cloverinfotech.com"
$ws.Range("D3").Value = "'7368768761"
$ws.Range("E3").Value = 9999

$ws.Range("A4").Value = "'3"
$ws.Range("B4").Value = "viraj"
$ws.Range("C4").Value = "viraj@gmail.com"
$ws.Range("D4").Value = "'8768768761"
$ws.Range("E4").Value = 190000

$ws.Range("A5").Value = "'4"
$ws.Range("B5").Value = "rajesh"
$ws.Range("C5").Value = "rajesh@gmail.com"
$ws.Range("D5").Value = "'83987987987"
$ws.Range("E5").Value = 10000000

$ws.Range("A6").Value = "'5"
$ws.Range("B6").Value = "bipin tiwariii"
$ws.Range("C6").Value = "bipin@gmail.com"
$ws.Range("D6").Value = "'82739879898"
$ws.Range("E6").Value = 9999999

$ws.Range("A7").Value = "'6"
$ws.Range("B7").Value = "saideep"
$ws.Range("C7").Value = "saideep@gmail.com"
$ws.Range("D7").Value = "'21837987981"
$ws.Range("E7").Value = 440000

$ws.Range("A8").Value = "'7"
$ws.Range("B8").Value = "test"
$ws.Range("C8").Value = "test@gmail.com"
$ws.Range("D8").Value = "'3009739711"
$ws.Range("E8").Value = 19991

# Remove the quote-prefix style Excel applied for the leading apostrophe so
# the cells keep the workbook's default style (index 0), matching the
# original text cells in the sheet.
$ws.Range("A2:A8").ClearFormats()
$ws.Range("D2:D8").ClearFormats()

# Resize the table to include the new rows
$tbl.Resize($ws.Range("A1:E8"))

Write-Host "Done"
Write-Host $tbl.Range.Address()
